$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel's
# General-format type inference (e.g. "1.01", "4.80") are forced to remain
# literal text by temporarily switching the cell to Text format, assigning the
# value, then clearing the format again so the cell keeps its original
# (unstyled / General) appearance while the stored value stays a string.

$ws.Range("D2").Value = "60.525.94"
$ws.Range("E2").Value = "  +6.38%  "
$ws.Range("D3").Value = "2.621.72"
$ws.Range("E3").Value = "  +8.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.10"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "2.657.86"
$ws.Range("E9").Value = "  +9.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.47"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.345"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "3.124.31"
$ws.Range("E14").Value = "  +10.46%  "
$ws.Range("D15").Value = "60.794.02"
$ws.Range("E15").Value = "  +6.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000141"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.18%  "
$ws.Range("D18").Value = "2.669.25"
$ws.Range("E18").Value = "  +10.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.57"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.425"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.30%  "
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "0.0₃0871"
$ws.Range("E28").Value = "  +11.73%  "
$ws.Range("E29").Value = "  +5.57%  "
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.58"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("E33").Value = "  +4.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.83"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.08"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.22"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "315.49"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +14.44%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.50"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.86%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.859"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.841"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +30.26%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.78"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +7.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.39"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.88%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.637"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.72%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0579"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.995"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.69"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +12.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.91"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.53%  "
$ws.Range("D49").Value = "2.064.61"
$ws.Range("E49").Value = "  +9.24%  "
$ws.Range("E50").Value = "  +4.43%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.30"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.98%  "
